# "drivers with no other expanded"
# - B13 was a placeholder "?" classification; resolve it to the proper
#   "frazzle (land ice) change" driver category (row's note explains it
#   needed to know what kind of ice to code, and A13 = "ice cover duration").
# - Rows whose Notes column says "not a category" previously left the
#   driver (column B) blank; make that explicit with "NA".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("B13").Value = "frazzle (land ice) change"

$ws1.Range("B17").Value = "NA"
$ws1.Range("B18").Value = "NA"
$ws1.Range("B19").Value = "NA"
$ws1.Range("B22").Value = "NA"
$ws1.Range("B23").Value = "NA"

# Leave the workbook with Sheet1 active and the cursor on the last edited
# cell, matching the saved view state in the edited file.
$ws1.Activate()
$ws1.Range("B23").Select() | Out-Null
